$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '27.187.52'
$ws.Range('D3').Value = '1.686.27'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '216.02'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('E6').Value = '  +0.70%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '23.14'
$ws.Range('E8').Value = '  +8.22%  '
$ws.Range('E9').Value = '  +3.30%  '
$ws.Range('E10').Value = '  +0.89%  '
$ws.Range('D11').Value = '0.0889'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = '1.924.07'
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').Value = '1.702.13'
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('D14').Value = '4.20'
$ws.Range('E14').Value = '  +2.36%  '
$ws.Range('E15').Value = '  +3.92%  '
$ws.Range('D16').Value = '66.92'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '27.192.03'
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('D18').Value = '236.24'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').Value = '8.04'
$ws.Range('E19').Value = '  -2.03%  '
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').Value = '4.57'
$ws.Range('E22').Value = '  +2.43%  '
$ws.Range('D23').Value = '9.62'
$ws.Range('E23').Value = '  +4.04%  '
$ws.Range('E24').Value = '  -3.12%  '
$ws.Range('D25').Value = '147.55'
$ws.Range('E25').Value = '  +0.48%  '
$ws.Range('D26').Value = '7.34'
$ws.Range('E26').Value = '  +1.38%  '
$ws.Range('E27').Value = '  +2.49%  '
$ws.Range('E28').Value = '  +0.69%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('E32').Value = '  +1.28%  '
$ws.Range('D33').Value = '1.544.00'
$ws.Range('E33').Value = '  +1.95%  '
$ws.Range('D34').Value = '3.25'
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('E36').Value = '  +2.63%  '
$ws.Range('E37').Value = '  +2.75%  '
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('E39').Value = '  -0.28%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '5.75'
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '69.12'
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('E44').Value = '  -1.02%  '
$ws.Range('D45').Value = '1.832.40'
$ws.Range('E45').Value = '  +0.62%  '
$ws.Range('D46').Value = '0.788'
$ws.Range('E46').Value = '  +0.73%  '
$ws.Range('D47').Value = '90.19'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('E48').Value = '  +4.91%  '
$ws.Range('E49').Value = '  +6.07%  '
$ws.Range('D50').Value = '8.26'
$ws.Range('E50').Value = '  +5.08%  '
$ws.Range('E51').Value = '  -0.78%  '
